# Refresh "Price" (D) and "Volume(1h)" (E) columns on the crypto symbol sheet.
# Values are stored as literal text (inlineStr) in the workbook, e.g. "261.10"
# and "1.57%", so each target cell is forced to Text format before the write
# and restored to the "Normal" style afterwards - this stops Excel/COM from
# re-interpreting the numeric-looking strings as real numbers/percentages.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2: BNB
Set-TextValue $ws.Range("D2") "261.10"
Set-TextValue $ws.Range("E2") "1.57%"

# Row 3: OKB
Set-TextValue $ws.Range("D3") "27.32"
Set-TextValue $ws.Range("E3") "1.09%"

# Row 4: HuobiToken
Set-TextValue $ws.Range("D4") "4.712"
Set-TextValue $ws.Range("E4") "4.30%"

# Row 5: Cronos
Set-TextValue $ws.Range("D5") "0.06080"
Set-TextValue $ws.Range("E5") "3.13%"

# Row 6: KuCoinToken
Set-TextValue $ws.Range("D6") "6.674"
Set-TextValue $ws.Range("E6") "0.96%"

# Row 7: MXToken
Set-TextValue $ws.Range("D7") "0.8458"
Set-TextValue $ws.Range("E7") "-0.56%"

# Row 8: FTXToken
Set-TextValue $ws.Range("D8") "0.9232"
Set-TextValue $ws.Range("E8") "-0.42%"

# Row 10: LiechtensteinCryptoassetsExchange
Set-TextValue $ws.Range("D10") "0.05038"
Set-TextValue $ws.Range("E10") "19.46%"

# Row 11: MandalaExchangeToken
Set-TextValue $ws.Range("D11") "0.07104"
Set-TextValue $ws.Range("E11") "1.47%"

# Row 12: BitrueCoin
Set-TextValue $ws.Range("D12") "0.03130"
Set-TextValue $ws.Range("E12") "2.62%"

# Row 13: BitMartToken
Set-TextValue $ws.Range("D13") "0.09072"
Set-TextValue $ws.Range("E13") "-0.30%"

# Row 14: BitForexToken
Set-TextValue $ws.Range("E14") "-0.21%"

# Row 15: One
Set-TextValue $ws.Range("D15") "0.0006081"
Set-TextValue $ws.Range("E15") "0.87%"

# Row 16: TigerCash
Set-TextValue $ws.Range("D16") "0.006121"
Set-TextValue $ws.Range("E16") "1.69%"

# Row 17: LEO
Set-TextValue $ws.Range("E17") "-0.48%"

# Row 18: GateToken
Set-TextValue $ws.Range("E18") "-0.73%"

# Row 19: BTSEToken
Set-TextValue $ws.Range("E19") "-2.02%"

# Row 20: BitpandaEcosystemToken
Set-TextValue $ws.Range("D20") "0.3127"
Set-TextValue $ws.Range("E20") "1.52%"

# Row 22: MCDex
Set-TextValue $ws.Range("D22") "4.090"
Set-TextValue $ws.Range("E22") "4.79%"

# Row 23: CoinExToken
Set-TextValue $ws.Range("D23") "0.04233"
Set-TextValue $ws.Range("E23") "-0.65%"

# Row 24: BitKan
Set-TextValue $ws.Range("D24") "0.001222"
Set-TextValue $ws.Range("E24") "0.34%"

# Row 25: HotbitToken
Set-TextValue $ws.Range("E25") "-9.06%"

# Row 26: NitroEx
Set-TextValue $ws.Range("E26") "0.10%"

# Row 27: UpBots
Set-TextValue $ws.Range("D27") "0.0001575"
Set-TextValue $ws.Range("E27") "3.44%"

# Row 40: IDEX
Set-TextValue $ws.Range("D40") "0.03868"
Set-TextValue $ws.Range("E40") "1.81%"

# Row 41: BKEXToken
Set-TextValue $ws.Range("D41") "0.1114"
Set-TextValue $ws.Range("E41") "1.24%"

# Row 42: KickToken
Set-TextValue $ws.Range("D42") "0.004095"
Set-TextValue $ws.Range("E42") "5.19%"

# Row 43: LocalTraders
Set-TextValue $ws.Range("E43") "21.58%"

# Row 44: CEJI
Set-TextValue $ws.Range("D44") "0.002221"
Set-TextValue $ws.Range("E44") "-7.02%"

# Row 45: CoinLion
Set-TextValue $ws.Range("D45") "0.00005277"
Set-TextValue $ws.Range("E45") "-1.33%"

# Row 46: Kangarootoken
Set-TextValue $ws.Range("E46") "0.08%"

# Row 47: CoinbaseStockToken
Set-TextValue $ws.Range("E47") "6.93%"

# Row 48: BOLO
Set-TextValue $ws.Range("D48") "0.1353"
Set-TextValue $ws.Range("E48") "-46.44%"

# Row 49: CryptobidCoin
Set-TextValue $ws.Range("E49") "0.08%"

# Row 50: SpecialPowerGold
Set-TextValue $ws.Range("E50") "0.08%"
